# "add logic to automate the date" - refresh the recorded run timestamps
# (Start/End/Duration) on the Summary sheet and for each TestCase row on
# the results sheet, as produced by a fresh automated test run.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$results = $wb.Worksheets.Item(2)

# Summary sheet: overall Start Time / End Time / Duration
$summary.Range("B6").Value = "2021-06-03T01:27:23 IST"
$summary.Range("B7").Value = "2021-06-03T01:28:21 IST"
$summary.Range("B8").Value = "57702 ms"

# Results sheet: per test-case Start Time / End Time / Duration
# Row 2 - loginPagetest
$results.Range("D2").Value = "2021-06-03T01:27:42 IST"
$results.Range("E2").Value = "2021-06-03T01:27:49 IST"
$results.Range("F2").Value = "7848 ms"

# Row 3 - logoutPageTest (starts right when the login test ends)
$results.Range("D3").Value = "2021-06-03T01:27:49 IST"
$results.Range("E3").Value = "2021-06-03T01:28:07 IST"
$results.Range("F3").Value = "17768 ms"
